$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet is being repurposed from "Create Item Category" to "Delete Item"
$ws.Name = "Delete Item"

# Drop any pre-existing formatting on the new header cells (B1:E1) so the
# row ends up with a single, consistent style; A1 already carries the
# correct header look from the template and is left untouched.
$ws.Range("B1:E1").ClearFormats()

# New header row: Item Code | Item Name | Category | Group | Unit
$ws.Range("A1").Value = "Item Code"
$ws.Range("B1").Value = "Item Name"
$ws.Range("C1").Value = "Category"
$ws.Range("D1").Value = "Group"
$ws.Range("E1").Value = "Unit"

# B1:E1 share the same fill as A1, without the font color override
$ws.Range("B1:E1").Interior.ThemeColor = 4

# Reset the selection back to the top-left cell (new default view state)
$ws.Range("A1").Select()

Write-Output "done"
